$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 0.092045493423938751
$ws.Range("C13").Value = 0.24352946877479553
$ws.Range("B14").Value = 0.27538201212882996
$ws.Range("C14").Value = 0.72440314292907715
$ws.Range("B15").Value = 0.89337277412414551
$ws.Range("C15").Value = 1.1186583042144775
$ws.Range("B16").Value = 2.6156282424926758
$ws.Range("C16").Value = 0.92864501476287842
$ws.Range("B17").Value = 4.997530460357666
$ws.Range("C17").Value = 1.0389989614486694
$ws.Range("B18").Value = 7.1618566513061523
$ws.Range("C18").Value = 2.0097074508666992
$ws.Range("B19").Value = 8.315032958984375
$ws.Range("C19").Value = 2.8133106231689453
$ws.Range("B20").Value = 8.0716428756713867
$ws.Range("C20").Value = 3.0895397663116455
$ws.Range("B21").Value = 6.575164794921875
$ws.Range("C21").Value = 2.9347128868103027
$ws.Range("B22").Value = 4.4180612564086914
$ws.Range("C22").Value = 2.5728747844696045
$ws.Range("B23").Value = 2.3786113262176514
$ws.Range("C23").Value = 1.9944329261779785
$ws.Range("B24").Value = 1.0038083791732788
$ws.Range("C24").Value = 1.1548430919647217
$ws.Range("B25").Value = 0.27529314160346985
$ws.Range("C25").Value = 0.433414101600647
$ws.Range("B26").Value = 0.017384422942996025
$ws.Range("C26").Value = 0.045994862914085388
$ws.Range("B30").Value = 0.024346509948372841
$ws.Range("C30").Value = 0.064414806663990021
$ws.Range("B31").Value = 0.064751468598842621
$ws.Range("C31").Value = 0.17131628096103668
$ws.Range("B32").Value = 0.11376496404409409
$ws.Range("C32").Value = 0.25475242733955383
$ws.Range("B33").Value = 0.14818465709686279
$ws.Range("C33").Value = 0.278164803981781
$ws.Range("B34").Value = 0.15036332607269287
$ws.Range("C34").Value = 0.25680005550384521
$ws.Range("B35").Value = 0.1225612685084343
$ws.Range("C35").Value = 0.23948028683662415
$ws.Range("B36").Value = 0.085873141884803772
$ws.Range("C36").Value = 0.22322170436382294
$ws.Range("B37").Value = 0.0701429694890976
$ws.Range("C37").Value = 0.15945154428482056
$ws.Range("B38").Value = 0.063696332275867462
$ws.Range("C38").Value = 0.1088627502322197
$ws.Range("B39").Value = 0.062066059559583664
$ws.Range("C39").Value = 0.14194025099277496
$ws.Range("B40").Value = 0.079142086207866669
$ws.Range("C40").Value = 0.14739426970481873
$ws.Range("B41").Value = 0.1279418021440506
$ws.Range("C41").Value = 0.16334426403045654
$ws.Range("B42").Value = 0.1684214174747467
$ws.Range("C42").Value = 0.27095073461532593
$ws.Range("B43").Value = 0.17201359570026398
$ws.Range("C43").Value = 0.32843098044395447
$ws.Range("B44").Value = 0.13324034214019775
$ws.Range("C44").Value = 0.27869126200675964
$ws.Range("B45").Value = 0.066950529813766479
$ws.Range("C45").Value = 0.16580019891262054
$ws.Range("B46").Value = 0.019646428525447845
$ws.Range("C46").Value = 0.051979564130306244
$ws.Range("B47").Value = 0.0035715179983526468
$ws.Range("C47").Value = 0.0094493487849831581
$ws.Range("B48").Value = 0.028116539120674133
$ws.Range("C48").Value = 0.074389368295669556
$ws.Range("B49").Value = 0.080853044986724854
$ws.Range("C49").Value = 0.17834228277206421
$ws.Range("B50").Value = 0.14877079427242279
$ws.Range("C50").Value = 0.28770923614501953
$ws.Range("B51").Value = 0.19853869080543518
$ws.Range("C51").Value = 0.36019706726074219
$ws.Range("B52").Value = 0.23839394748210907
$ws.Range("C52").Value = 0.347348690032959
$ws.Range("B53").Value = 0.30619516968727112
$ws.Range("C53").Value = 0.41468191146850586
$ws.Range("B54").Value = 0.35998058319091797
$ws.Range("C54").Value = 0.69348597526550293
$ws.Range("B55").Value = 0.37021887302398682
$ws.Range("C55").Value = 0.86900711059570313
$ws.Range("B56").Value = 0.35516220331192017
$ws.Range("C56").Value = 0.78357130289077759
$ws.Range("B57").Value = 0.33036574721336365
$ws.Range("C57").Value = 0.54024893045425415
$ws.Range("B58").Value = 0.41762274503707886
$ws.Range("C58").Value = 0.46683052182197571
$ws.Range("B59").Value = 0.632537305355072
$ws.Range("C59").Value = 0.86840558052062988
$ws.Range("B60").Value = 0.8666996955871582
$ws.Range("C60").Value = 1.3545219898223877
$ws.Range("B61").Value = 0.94491147994995117
$ws.Range("C61").Value = 1.681327223777771
$ws.Range("B62").Value = 0.83280318975448608
$ws.Range("C62").Value = 1.6705116033554077
$ws.Range("B63").Value = 0.59954154491424561
$ws.Range("C63").Value = 1.296189546585083
$ws.Range("B64").Value = 0.39407306909561157
$ws.Range("C64").Value = 0.72886854410171509
$ws.Range("B65").Value = 0.29769939184188843
$ws.Range("C65").Value = 0.56754404306411743
$ws.Range("B66").Value = 0.3408207893371582
$ws.Range("C66").Value = 0.90172708034515381
$ws.Range("B67").Value = 0.43646791577339172
$ws.Range("C67").Value = 1.1426857709884644
$ws.Range("B68").Value = 0.73668259382247925
$ws.Range("C68").Value = 1.3019261360168457
$ws.Range("B69").Value = 1.1771696805953979
$ws.Range("C69").Value = 2.1689536571502686
$ws.Range("B70").Value = 1.568972110748291
$ws.Range("C70").Value = 3.4113447666168213
$ws.Range("B71").Value = 1.7376065254211426
$ws.Range("C71").Value = 4.225677490234375
$ws.Range("B72").Value = 1.6324841976165771
$ws.Range("C72").Value = 4.15089225769043
